$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two missing values (D146/E146) for the last pre-existing row.
$ws.Range("D146").Value = 0
$ws.Range("E146").Value = 0

# Append new fly/pupae observation rows (vial, treatment, time_hours, males, females)
# for time_hours 341 and 356 cohorts (rows 147-204).

$data = @(
    ,@(147, 1, "conditioned", 341, 6, 6)
    ,@(148, 1, "unconditioned", 341, 0, 0)
    ,@(149, 2, "conditioned", 341, 3, 8)
    ,@(150, 2, "unconditioned", 341, 0, 0)
    ,@(151, 3, "conditioned", 341, 3, 4)
    ,@(152, 3, "unconditioned", 341, 0, 0)
    ,@(153, 4, "conditioned", 341, 1, 2)
    ,@(154, 4, "unconditioned", 341, 3, 4)
    ,@(155, 5, "conditioned", 341, 0, 0)
    ,@(156, 5, "unconditioned", 341, 5, 6)
    ,@(157, 6, "conditioned", 341, 5, 6)
    ,@(158, 6, "unconditioned", 341, 2, 0)
    ,@(159, 7, "conditioned", 341, 0, 0)
    ,@(160, 7, "unconditioned", 341, 7, 3)
    ,@(161, 8, "conditioned", 341, 4, 7)
    ,@(162, 8, "unconditioned", 341, 6, 4)
    ,@(163, 9, "conditioned", 341, 2, 2)
    ,@(164, 9, "unconditioned", 341, 5, 2)
    ,@(165, 10, "conditioned", 341, 1, 4)
    ,@(166, 10, "unconditioned", 341, 3, 4)
    ,@(167, 11, "conditioned", 341, 3, 1)
    ,@(168, 11, "unconditioned", 341, 4, 1)
    ,@(169, 12, "conditioned", 341, 1, 2)
    ,@(170, 12, "unconditioned", 341, 6, 2)
    ,@(171, 13, "conditioned", 341, 0, 0)
    ,@(172, 13, "unconditioned", 341, 1, 4)
    ,@(173, 14, "conditioned", 341, 4, 3)
    ,@(174, 14, "unconditioned", 341, 0, 0)
    ,@(175, 15, "unconditioned", 341, 1, 0)
    ,@(176, 1, "conditioned", 356, 0, 2)
    ,@(177, 1, "unconditioned", 356, 0, 0)
    ,@(178, 2, "conditioned", 356, 0, 3)
    ,@(179, 2, "unconditioned", 356, 0, 0)
    ,@(180, 3, "conditioned", 356, 1, 1)
    ,@(181, 3, "unconditioned", 356, 0, 0)
    ,@(182, 4, "conditioned", 356, 4, 4)
    ,@(183, 4, "unconditioned", 356, 2, 4)
    ,@(184, 5, "conditioned", 356, 1, 0)
    ,@(185, 5, "unconditioned", 356, 3, 0)
    ,@(186, 6, "conditioned", 356, 1, 1)
    ,@(187, 6, "unconditioned", 356, 4, 2)
    ,@(188, 7, "conditioned", 356, 0, 0)
    ,@(189, 7, "unconditioned", 356, 4, 2)
    ,@(190, 8, "conditioned", 356, 3, 1)
    ,@(191, 8, "unconditioned", 356, 0, 1)
    ,@(192, 9, "conditioned", 356, 1, 1)
    ,@(193, 9, "unconditioned", 356, 1, 0)
    ,@(194, 10, "conditioned", 356, 0, 2)
    ,@(195, 10, "unconditioned", 356, 2, 2)
    ,@(196, 11, "conditioned", 356, 1, 1)
    ,@(197, 11, "unconditioned", 356, 0, 0)
    ,@(198, 12, "conditioned", 356, 0, 2)
    ,@(199, 12, "unconditioned", 356, 1, 0)
    ,@(200, 13, "conditioned", 356, 0, 0)
    ,@(201, 13, "unconditioned", 356, 1, 1)
    ,@(202, 14, "conditioned", 356, 2, 2)
    ,@(203, 14, "unconditioned", 356, 0, 0)
    ,@(204, 15, "unconditioned", 356, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Match the final selection/active cell reported after the edit.
$ws.Range("F204").Select()
